$d = $word.ActiveDocument

function Format-ChangedRun {
    param(
        [string]$Context,
        [int]$RunLen,
        [bool]$MakeBold
    )

    $rng = $d.Content
    $found = $rng.Find.Execute($Context, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Output "NOT FOUND: $Context"
        return
    }

    $runStart = $rng.End - $RunLen
    $runRange = $d.Range($runStart, $rng.End)

    # Clear the explicit "strike-through = off" direct formatting (matches the
    # author's later resave, which drops the redundant <w:strike w:val="0"/>).
    $runRange.Font.StrikeThrough = $false

    if ($MakeBold) {
        $runRange.Font.Bold = $true
    }

    Write-Output "Formatted [$($runRange.Text)] bold=$MakeBold"
}

# The heading run " changed" (after "What is LOREM IPSUM") also picks up bold,
# matching the rest of the bold heading text around it.
Format-ChangedRun "LOREM IPSUM changed" 8 $true

# The remaining five "changed here[ too]" runs only lose the redundant
# strike-through override; their bold state is untouched.
Format-ChangedRun "Letraset sheets containing Lorem Ipsum changed here too" 17 $false
Format-ChangedRun "including versions of Lorem Ipsum changed here" 13 $false
Format-ChangedRun "default model changed here too" 17 $false
Format-ChangedRun "obscure Latin words changed here too" 17 $false
Format-ChangedRun "Bonorum et changed here too Malorum" 17 $false
